$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 9288.299999999999, 10036, 8426, 0.191639502843221),
    @(1, 9425.533333333333, 10063, 8369, 0.2042194366455078),
    @(2, 9568, 10247, 8824, 0.2024188915888468),
    @(3, 9873.133333333333, 10544, 8754, 0.1980333248774211),
    @(4, 9114.1, 10237, 7908, 0.1903159856796265),
    @(5, 10040.93333333333, 10948, 9161, 0.2021841923395793),
    @(6, 9269.666666666666, 9893, 8459, 0.2013073523839315),
    @(7, 9117.966666666667, 10194, 8317, 0.1996895790100098),
    @(8, 8664.933333333332, 9498, 7300, 0.1919730265935262),
    @(9, 8893.633333333333, 9790, 8124, 0.1937473853429159)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
